# "small upd + pthw"
#
# Updates homework text for several rows, marks row 11 (Литература) as a
# rest day, expands on the "metal from ore" task for row 14 (Химия),
# tidies the spacing in the alexlarin.net link on row 5 (Мат. анализ),
# widens column B to fit the longer text, tweaks the zoom level and moves
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Мат. анализ (row 5): tidy up spacing inside the parentheses of the link.
$ws.Range("B5").Value = "см. почту + сделать 267 вариант (http://alexlarin.net/ege/2019/trvar267.html)"

# Литература (row 11): nothing assigned today - everyone rests.
$ws.Range("B11").Value = "Отдыхаем"

# Химия (row 14): clarify that only students who haven't passed need to prep.
$ws.Range("B14").Value = "Готовиться к устной защите способов получение металла из руды ( тем, кто не сдал )"

# Column B needs to be a bit wider to comfortably fit the updated text.
$ws.Columns.Item(2).ColumnWidth = 98.17

# Slightly reduce the zoom level and move the active selection/cursor.
$excel.ActiveWindow.Zoom = 105
$ws.Range("B15").Select() | Out-Null
